$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Sexo" (Gender) labels used in column G.
# Previously: Masculino / Femenino
# Now:        Mujer / Hombre  (with the underlying data swapped so that
# rows that used to read "Masculino" now read "Hombre" and rows that used
# to read "Femenino" now read "Mujer")
$lastRow = $ws.UsedRange.Rows.Count()
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G ("Sexo")
    $v = $cell.Value()
    if ($v -eq "Masculino") {
        $cell.Value = "Hombre"
    } elseif ($v -eq "Femenino") {
        $cell.Value = "Mujer"
    }
}

# Update the selected range / active cell to match the saved view state.
$ws.Range("S16").Select()
